# Edit: "How to make an app for Zendesk API v 2" -> "How to make an app for Zendesk"
# (title text is re-typed as 3 runs) plus a date-placeholder refresh
# (4/22/2013 -> 4/30/2013) across the slide master and every slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 1 title: drop the trailing " API v 2" and re-enter "for " /
#    "Zendesk" as their own runs (mirrors retyping the end of the title).
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

$prefix = "How to make an app "
$mid    = "for "
$last   = "Zendesk"
$wanted = $prefix + $mid + $last

$current = $titleRange.Text
if ($current.StartsWith($wanted)) {
    $extraLen = $current.Length - $wanted.Length
    if ($extraLen -gt 0) {
        $extraRange = $titleRange.Characters($wanted.Length + 1, $extraLen)
        $extraRange.Text = ""
    }

    $midRange = $titleRange.Characters($prefix.Length + 1, $mid.Length)
    $midRange.Text = $mid

    $lastRange = $titleRange.Characters($prefix.Length + $mid.Length + 1, $last.Length)
    $lastRange.Text = $last
}

# ---------------------------------------------------------------------------
# 2) Refresh the cached "datetimeFigureOut" date placeholder text from
#    4/22/2013 to 4/30/2013 on the slide master and on every slide layout.
# ---------------------------------------------------------------------------
function Get-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shape = $container.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            try {
                if ($shape.PlaceholderFormat.Type -eq 16) {
                    return $shape
                }
            } catch {
                # not a placeholder shape - ignore
            }
        }
    }
    return $null
}

$newDate = "4/30/2013"
$master = $p.SlideMaster

$masterDateShape = Get-DatePlaceholder($master)
if ($masterDateShape -ne $null) {
    if ($masterDateShape.TextFrame.TextRange.Text -eq "4/22/2013") {
        $masterDateShape.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    $layoutDateShape = Get-DatePlaceholder($layout)
    if ($layoutDateShape -ne $null) {
        if ($layoutDateShape.TextFrame.TextRange.Text -eq "4/22/2013") {
            $layoutDateShape.TextFrame.TextRange.Text = $newDate
        }
    }
}
